$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("M2").Value = 1.847479
$ws.Range("N2").Value = 5.542437
$ws.Range("O2").Value = 0.3181373042830636
$ws.Range("P2").Value = 0.3181373042830637
$ws.Range("Q2").Value = 11.46548300201133
$ws.Range("R2").Value = 103.189347018102
$ws.Range("S2").Value = 0.004777544383605788
$ws.Range("T2").Value = 0.004777544383605788
$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("N3").Value = 6.300930999999999
$ws.Range("O3").Value = 0.3616750542791174
$ws.Range("P3").Value = 0.3616750542791174
$ws.Range("Q3").Value = 13.03455813342511
$ws.Range("R3").Value = 117.311023200826
$ws.Range("S3").Value = 0.005431361242453022
$ws.Range("T3").Value = 0.005431361242453021
$ws.Range("G4").Value = 6.206015333333333
$ws.Range("H4").Value = 18.618046
$ws.Range("I4").Value = 0.0150172404156507
$ws.Range("J4").Value = 0.0150172404156507
$ws.Range("M4").Value = 1.859385666666666
$ws.Range("N4").Value = 5.578156999999999
$ws.Range("O4").Value = 0.320187641437819
$ws.Range("P4").Value = 0.320187641437819
$ws.Range("Q4").Value = 11.53937595791355
$ws.Range("R4").Value = 103.854383621222
$ws.Range("S4").Value = 0.00480833478959189
$ws.Range("T4").Value = 0.00480833478959189
$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("M5").Value = 1.847479
$ws.Range("N5").Value = 5.542437
$ws.Range("O5").Value = 0.3181373042830636
$ws.Range("P5").Value = 0.3181373042830637
$ws.Range("Q5").Value = 711.3763528389563
$ws.Range("R5").Value = 6402.387175550607
$ws.Range("S5").Value = 0.2964229329492286
$ws.Range("T5").Value = 0.2964229329492286
$ws.Range("I6").Value = 0.9317452840597572
$ws.Range("J6").Value = 0.9317452840597571
$ws.Range("N6").Value = 6.300930999999999
$ws.Range("O6").Value = 0.3616750542791174
$ws.Range("P6").Value = 0.3616750542791174
$ws.Range("R6").Value = 7278.567140849641
$ws.Range("S6").Value = 0.3369890261866244
$ws.Range("T6").Value = 0.3369890261866243
$ws.Range("I7").Value = 0.9317452840597572
$ws.Range("J7").Value = 0.9317452840597571
$ws.Range("M7").Value = 1.859385666666666
$ws.Range("N7").Value = 5.578156999999999
$ws.Range("O7").Value = 0.320187641437819
$ws.Range("P7").Value = 0.320187641437819
$ws.Range("Q7").Value = 715.9610442523918
$ws.Range("R7").Value = 6443.649398271526
$ws.Range("S7").Value = 0.2983333249239044
$ws.Range("T7").Value = 0.2983333249239044
$ws.Range("G8").Value = 22.00088566666667
$ws.Range("H8").Value = 66.002657
$ws.Range("I8").Value = 0.05323747552459213
$ws.Range("J8").Value = 0.05323747552459213
$ws.Range("M8").Value = 1.847479
$ws.Range("N8").Value = 5.542437
$ws.Range("O8").Value = 0.3181373042830636
$ws.Range("P8").Value = 0.3181373042830637
$ws.Range("Q8").Value = 40.64617425056766
$ws.Range("R8").Value = 365.815568255109
$ws.Range("S8").Value = 0.01693682695022932
$ws.Range("T8").Value = 0.01693682695022932
$ws.Range("G9").Value = 22.00088566666667
$ws.Range("H9").Value = 66.002657
$ws.Range("I9").Value = 0.05323747552459213
$ws.Range("J9").Value = 0.05323747552459213
$ws.Range("N9").Value = 6.300930999999999
$ws.Range("O9").Value = 0.3616750542791174
$ws.Range("P9").Value = 0.3616750542791174
$ws.Range("Q9").Value = 46.20868750818522
$ws.Range("R9").Value = 415.8781875736669
$ws.Range("S9").Value = 0.01925466685004004
$ws.Range("T9").Value = 0.01925466685004004
$ws.Range("G10").Value = 22.00088566666667
$ws.Range("H10").Value = 66.002657
$ws.Range("I10").Value = 0.05323747552459213
$ws.Range("J10").Value = 0.05323747552459213
$ws.Range("M10").Value = 1.859385666666666
$ws.Range("N10").Value = 5.578156999999999
$ws.Range("O10").Value = 0.320187641437819
$ws.Range("P10").Value = 0.320187641437819
$ws.Range("Q10").Value = 40.9081314625721
$ws.Range("R10").Value = 368.173183163149
$ws.Range("S10").Value = 0.01704598172432277
$ws.Range("T10").Value = 0.01704598172432277
